$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Insert a brand-new "2022-Q4" sheet right before the existing "2022-Q3"
#    sheet. Duplicating "2022-Q3" first gives us an exact copy of all the
#    shared formatting (header style, border style, column A numbering
#    style, page margins, etc.) which we then overwrite with the 2022-Q4
#    numbers.
# ---------------------------------------------------------------------------
$q3 = $wb.Worksheets.Item("2022-Q3")
$q3.Copy($q3)
$newSheet = $wb.Worksheets.Item("2022-Q3 (2)")
$newSheet.Name = "2022-Q4"

# The copied sheet has 18 data rows (rows 2-19); 2022-Q4 only needs 8 data
# rows (rows 2-9), so drop the extra rows entirely.
$newSheet.Range("A10:H19").Clear()

# Columns B-G hold text (fund code / name / percentages formatted as text,
# matching every other quarter sheet in this workbook) - force text storage
# so things like leading zeros ("014007") and trailing zeros ("93.81",
# "0.4579" etc.) round-trip exactly instead of being parsed as numbers.
$newSheet.Range("B2:G9").NumberFormat = "@"

$newRows = @(
    @(0, "014007", "华安制造升级一年持有混合A",       "11.68", "93.81", "3.92", "0.4579", 8),
    @(1, "040001", "华安创新混合",                     "14.64", "74.94", "3.03", "0.4436", 8),
    @(2, "006154", "华安制造先锋混合A",                 "10.24", "93.81", "3.92", "0.4014", 10),
    @(3, "160425", "华安创业板两年定期开放混合",         "1.71",  "97.00", "4.65", "0.0795", 8),
    @(4, "014008", "华安制造升级一年持有混合C",          "0.55",  "93.81", "3.92", "0.0216", 8),
    @(5, "013507", "华安制造先锋混合C",                  "0.37",  "93.81", "3.92", "0.0145", 10),
    @(6, "000531", "东吴阿尔法灵活配置混合A",            "0.26",  "84.60", "3.49", "0.0091", 9),
    @(7, "014581", "东吴阿尔法灵活配置混合C",            "0.03",  "84.60", "3.49", "0.0010", 9)
)

$r = 2
foreach ($row in $newRows) {
    $newSheet.Cells.Item($r, 1).Value = $row[0]
    $newSheet.Cells.Item($r, 2).Value = $row[1]
    $newSheet.Cells.Item($r, 3).Value = $row[2]
    $newSheet.Cells.Item($r, 4).Value = $row[3]
    $newSheet.Cells.Item($r, 5).Value = $row[4]
    $newSheet.Cells.Item($r, 6).Value = $row[5]
    $newSheet.Cells.Item($r, 7).Value = $row[6]
    $newSheet.Cells.Item($r, 8).Value = $row[7]
    $r++
}

# Drop the text-format override now that the strings are safely stored, so
# the cells end up with the same "no explicit style" look the rest of the
# sheet has.
$newSheet.Range("B2:G9").ClearFormats()

# ---------------------------------------------------------------------------
# 2) Prepend a "2022-Q4" summary row to the "总计" (totals) sheet, shifting
#    the eight existing rows down by one.
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")

# Column A is just a plain 0-based row counter (row 2 -> 0, row 3 -> 1, ...),
# independent of the data, so only shift columns B:D down and re-number A
# afterwards instead of dragging the old index values down with the rest of
# the row.
for ($i = 9; $i -ge 2; $i--) {
    $srcRow = $summary.Range("B" + $i + ":D" + $i)
    $dstRow = $summary.Range("B" + ($i + 1) + ":D" + ($i + 1))
    $srcRow.Copy($dstRow)
}

$summary.Cells.Item(2, 2).Value = "2022-Q4"
$summary.Cells.Item(2, 3).Value = 8
$summary.Cells.Item(2, 4).Value = 1.43

# Row 10 is brand new, so A10 needs the same style as the rest of column A
# (A2:A9) - copy just the formatting down from A9.
$summary.Range("A9").Copy()
$summary.Range("A10").PasteSpecial(-4122)

for ($i = 2; $i -le 10; $i++) {
    $summary.Cells.Item($i, 1).Value = $i - 2
}

# ---------------------------------------------------------------------------
# 3) Restore the active-tab selection to "2020-Q4" (the last sheet), which
#    is where it was before the new sheet got inserted/copied-in.
# ---------------------------------------------------------------------------
$wb.Worksheets.Item("2020-Q4").Activate()
